# Generate Report for Handoff
# The d22d2421-... file has finished translation and is now ready for handoff.
# Update its status/priority/handoff-timestamp across the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 16:11:42"

# Columns E/F widen slightly to fit the new, longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 16.4
$wsOverview.Columns.Item(6).ColumnWidth = 16.4

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-18 16:11:38"

# Column C widens to fit the new, longer status text.
$wsZhCn.Columns.Item(3).ColumnWidth = 16.4

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-18 16:11:42"

# Column C widens to fit the new, longer status text.
$wsDeDe.Columns.Item(3).ColumnWidth = 16.4
